$d = $word.ActiveDocument

$pairs = @(
    @("82÷9=9, 1", "87÷2=43, 1"),
    @("26÷3=8, 2", "12÷3=4, 0"),
    @("35÷4=8, 3", "86÷2=43, 0"),
    @("32÷7=4, 4", "33÷5=6, 3"),
    @("54÷6=9, 0", "87÷2=43, 1"),
    @("62÷7=8, 6", "33÷3=11, 0"),
    @("58÷3=19, 1", "99÷2=49, 1"),
    @("68÷9=7, 5", "90÷6=15, 0"),
    @("95÷8=11, 7", "47÷4=11, 3"),
    @("73÷6=12, 1", "19÷7=2, 5"),
    @("91÷7=13, 0", "74÷2=37, 0"),
    @("12÷6=2, 0", "75÷7=10, 5"),
    @("64÷2=32, 0", "58÷6=9, 4"),
    @("45÷5=9, 0", "66÷6=11, 0"),
    @("41÷5=8, 1", "44÷7=6, 2"),
    @("89÷3=29, 2", "42÷9=4, 6"),
    @("53÷2=26, 1", "93÷2=46, 1"),
    @("47÷3=15, 2", "85÷4=21, 1"),
    @("48÷3=16, 0", "82÷7=11, 5"),
    @("47÷8=5, 7", "88÷4=22, 0"),
    @("51÷5=10, 1", "67÷5=13, 2"),
    @("85÷8=10, 5", "76÷6=12, 4"),
    @("10÷9=1, 1", "66÷5=13, 1"),
    @("19÷4=4, 3", "30÷5=6, 0"),
    @("64÷5=12, 4", "39÷2=19, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
